$d = $word.ActiveDocument
$d.Content.Find.Execute("61-15=46", $true, $false, $false, $false, $false, $false, 1, $false, "0+57=57", 2) | Out-Null
$d.Content.Find.Execute("11+47=58", $true, $false, $false, $false, $false, $false, 1, $false, "85-59=26", 2) | Out-Null
$d.Content.Find.Execute("21+59=80", $true, $false, $false, $false, $false, $false, 1, $false, "5+68=73", 2) | Out-Null
$d.Content.Find.Execute("85-49=36", $true, $false, $false, $false, $false, $false, 1, $false, "58-24=34", 2) | Out-Null
$d.Content.Find.Execute("46-5=41", $true, $false, $false, $false, $false, $false, 1, $false, "52+1=53", 2) | Out-Null
$d.Content.Find.Execute("49+50=99", $true, $false, $false, $false, $false, $false, 1, $false, "96-79=17", 2) | Out-Null
$d.Content.Find.Execute("95-15=80", $true, $false, $false, $false, $false, $false, 1, $false, "70-41=29", 2) | Out-Null
$d.Content.Find.Execute("13+44=57", $true, $false, $false, $false, $false, $false, 1, $false, "12+29=41", 2) | Out-Null
$d.Content.Find.Execute("89-56=33", $true, $false, $false, $false, $false, $false, 1, $false, "39-11=28", 2) | Out-Null
$d.Content.Find.Execute("6+22=28", $true, $false, $false, $false, $false, $false, 1, $false, "90-68=22", 2) | Out-Null
$d.Content.Find.Execute("51+17=68", $true, $false, $false, $false, $false, $false, 1, $false, "3+28=31", 2) | Out-Null
$d.Content.Find.Execute("15+37=52", $true, $false, $false, $false, $false, $false, 1, $false, "20+49=69", 2) | Out-Null
$d.Content.Find.Execute("64+19=83", $true, $false, $false, $false, $false, $false, 1, $false, "57+38=95", 2) | Out-Null
$d.Content.Find.Execute("96-18=78", $true, $false, $false, $false, $false, $false, 1, $false, "5+51=56", 2) | Out-Null
$d.Content.Find.Execute("76-33=43", $true, $false, $false, $false, $false, $false, 1, $false, "82-30=52", 2) | Out-Null
$d.Content.Find.Execute("74-54=20", $true, $false, $false, $false, $false, $false, 1, $false, "32+63=95", 2) | Out-Null
$d.Content.Find.Execute("60-33=27", $true, $false, $false, $false, $false, $false, 1, $false, "30+67=97", 2) | Out-Null
$d.Content.Find.Execute("82-5=77", $true, $false, $false, $false, $false, $false, 1, $false, "96-11=85", 2) | Out-Null
$d.Content.Find.Execute("38+3=41", $true, $false, $false, $false, $false, $false, 1, $false, "62-60=2", 2) | Out-Null
$d.Content.Find.Execute("22+66=88", $true, $false, $false, $false, $false, $false, 1, $false, "79-11=68", 2) | Out-Null
$d.Content.Find.Execute("62+26=88", $true, $false, $false, $false, $false, $false, 1, $false, "38+60=98", 2) | Out-Null
$d.Content.Find.Execute("76-65=11", $true, $false, $false, $false, $false, $false, 1, $false, "94-15=79", 2) | Out-Null
$d.Content.Find.Execute("26+58=84", $true, $false, $false, $false, $false, $false, 1, $false, "90-50=40", 2) | Out-Null
$d.Content.Find.Execute("18+13=31", $true, $false, $false, $false, $false, $false, 1, $false, "23+51=74", 2) | Out-Null
$d.Content.Find.Execute("34+38=72", $true, $false, $false, $false, $false, $false, 1, $false, "78+9=87", 2) | Out-Null
$d.Content.Find.Execute("73+0=73", $true, $false, $false, $false, $false, $false, 1, $false, "14+82=96", 2) | Out-Null
$d.Content.Find.Execute("1+85=86", $true, $false, $false, $false, $false, $false, 1, $false, "56+13=69", 2) | Out-Null
$d.Content.Find.Execute("79-20=59", $true, $false, $false, $false, $false, $false, 1, $false, "40-18=22", 2) | Out-Null
$d.Content.Find.Execute("42+2=44", $true, $false, $false, $false, $false, $false, 1, $false, "83-58=25", 2) | Out-Null
$d.Content.Find.Execute("33-6=27", $true, $false, $false, $false, $false, $false, 1, $false, "25+39=64", 2) | Out-Null
$d.Content.Find.Execute("54+37=91", $true, $false, $false, $false, $false, $false, 1, $false, "6+28=34", 2) | Out-Null
$d.Content.Find.Execute("94-76=18", $true, $false, $false, $false, $false, $false, 1, $false, "36-20=16", 2) | Out-Null
$d.Content.Find.Execute("69-39=30", $true, $false, $false, $false, $false, $false, 1, $false, "7+10=17", 2) | Out-Null
$d.Content.Find.Execute("62-45=17", $true, $false, $false, $false, $false, $false, 1, $false, "19+63=82", 2) | Out-Null
$d.Content.Find.Execute("85-63=22", $true, $false, $false, $false, $false, $false, 1, $false, "15+69=84", 2) | Out-Null
$d.Content.Find.Execute("80-36=44", $true, $false, $false, $false, $false, $false, 1, $false, "62-7=55", 2) | Out-Null
$d.Content.Find.Execute("40+22=62", $true, $false, $false, $false, $false, $false, 1, $false, "3+53=56", 2) | Out-Null
$d.Content.Find.Execute("70-65=5", $true, $false, $false, $false, $false, $false, 1, $false, "92-45=47", 2) | Out-Null
$d.Content.Find.Execute("6+48=54", $true, $false, $false, $false, $false, $false, 1, $false, "36+16=52", 2) | Out-Null
$d.Content.Find.Execute("3+39=42", $true, $false, $false, $false, $false, $false, 1, $false, "29+37=66", 2) | Out-Null
$d.Content.Find.Execute("48+20=68", $true, $false, $false, $false, $false, $false, 1, $false, "66-46=20", 2) | Out-Null
$d.Content.Find.Execute("3+2=5", $true, $false, $false, $false, $false, $false, 1, $false, "81-23=58", 2) | Out-Null
$d.Content.Find.Execute("70-42=28", $true, $false, $false, $false, $false, $false, 1, $false, "6+81=87", 2) | Out-Null
$d.Content.Find.Execute("84-82=2", $true, $false, $false, $false, $false, $false, 1, $false, "76+23=99", 2) | Out-Null
$d.Content.Find.Execute("51-21=30", $true, $false, $false, $false, $false, $false, 1, $false, "43+13=56", 2) | Out-Null
$d.Content.Find.Execute("76-75=1", $true, $false, $false, $false, $false, $false, 1, $false, "33+52=85", 2) | Out-Null
$d.Content.Find.Execute("5+1=6", $true, $false, $false, $false, $false, $false, 1, $false, "83-72=11", 2) | Out-Null
$d.Content.Find.Execute("45+14=59", $true, $false, $false, $false, $false, $false, 1, $false, "63-10=53", 2) | Out-Null
$d.Content.Find.Execute("13+65=78", $true, $false, $false, $false, $false, $false, 1, $false, "21-16=5", 2) | Out-Null
$d.Content.Find.Execute("13+80=93", $true, $false, $false, $false, $false, $false, 1, $false, "67-43=24", 2) | Out-Null
$d.Content.Find.Execute("66+5=71", $true, $false, $false, $false, $false, $false, 1, $false, "94-12=82", 2) | Out-Null
$d.Content.Find.Execute("17+23=40", $true, $false, $false, $false, $false, $false, 1, $false, "52-4=48", 2) | Out-Null
$d.Content.Find.Execute("7+5=12", $true, $false, $false, $false, $false, $false, 1, $false, "29-3=26", 2) | Out-Null
$d.Content.Find.Execute("36+41=77", $true, $false, $false, $false, $false, $false, 1, $false, "37+5=42", 2) | Out-Null
$d.Content.Find.Execute("44+10=54", $true, $false, $false, $false, $false, $false, 1, $false, "99-9=90", 2) | Out-Null
$d.Content.Find.Execute("49+22=71", $true, $false, $false, $false, $false, $false, 1, $false, "78-11=67", 2) | Out-Null
$d.Content.Find.Execute("40-31=9", $true, $false, $false, $false, $false, $false, 1, $false, "70-0=70", 2) | Out-Null
$d.Content.Find.Execute("11-6=5", $true, $false, $false, $false, $false, $false, 1, $false, "11+30=41", 2) | Out-Null
$d.Content.Find.Execute("89-40=49", $true, $false, $false, $false, $false, $false, 1, $false, "1+49=50", 2) | Out-Null
$d.Content.Find.Execute("64+6=70", $true, $false, $false, $false, $false, $false, 1, $false, "6-3=3", 2) | Out-Null
$d.Content.Find.Execute("19-17=2", $true, $false, $false, $false, $false, $false, 1, $false, "31+12=43", 2) | Out-Null
$d.Content.Find.Execute("78-37=41", $true, $false, $false, $false, $false, $false, 1, $false, "66-53=13", 2) | Out-Null
$d.Content.Find.Execute("16+19=35", $true, $false, $false, $false, $false, $false, 1, $false, "36-18=18", 2) | Out-Null
$d.Content.Find.Execute("47+33=80", $true, $false, $false, $false, $false, $false, 1, $false, "87-67=20", 2) | Out-Null
$d.Content.Find.Execute("58-53=5", $true, $false, $false, $false, $false, $false, 1, $false, "37-0=37", 2) | Out-Null
$d.Content.Find.Execute("96+1=97", $true, $false, $false, $false, $false, $false, 1, $false, "97-55=42", 2) | Out-Null
$d.Content.Find.Execute("30+29=59", $true, $false, $false, $false, $false, $false, 1, $false, "11+3=14", 2) | Out-Null
$d.Content.Find.Execute("59-8=51", $true, $false, $false, $false, $false, $false, 1, $false, "73-66=7", 2) | Out-Null
$d.Content.Find.Execute("68+11=79", $true, $false, $false, $false, $false, $false, 1, $false, "98-59=39", 2) | Out-Null
$d.Content.Find.Execute("38+14=52", $true, $false, $false, $false, $false, $false, 1, $false, "43-19=24", 2) | Out-Null
$d.Content.Find.Execute("37+6=43", $true, $false, $false, $false, $false, $false, 1, $false, "29+13=42", 2) | Out-Null
$d.Content.Find.Execute("91-43=48", $true, $false, $false, $false, $false, $false, 1, $false, "59-46=13", 2) | Out-Null
$d.Content.Find.Execute("54+21=75", $true, $false, $false, $false, $false, $false, 1, $false, "52-47=5", 2) | Out-Null
$d.Content.Find.Execute("54-16=38", $true, $false, $false, $false, $false, $false, 1, $false, "10+89=99", 2) | Out-Null
$d.Content.Find.Execute("26+7=33", $true, $false, $false, $false, $false, $false, 1, $false, "84-15=69", 2) | Out-Null
$d.Content.Find.Execute("55-30=25", $true, $false, $false, $false, $false, $false, 1, $false, "25+8=33", 2) | Out-Null
$d.Content.Find.Execute("35+47=82", $true, $false, $false, $false, $false, $false, 1, $false, "20+2=22", 2) | Out-Null
$d.Content.Find.Execute("82-47=35", $true, $false, $false, $false, $false, $false, 1, $false, "13+25=38", 2) | Out-Null
$d.Content.Find.Execute("70-8=62", $true, $false, $false, $false, $false, $false, 1, $false, "41+28=69", 2) | Out-Null
$d.Content.Find.Execute("44+52=96", $true, $false, $false, $false, $false, $false, 1, $false, "88+7=95", 2) | Out-Null
$d.Content.Find.Execute("90-64=26", $true, $false, $false, $false, $false, $false, 1, $false, "17+77=94", 2) | Out-Null
$d.Content.Find.Execute("71-58=13", $true, $false, $false, $false, $false, $false, 1, $false, "87-14=73", 2) | Out-Null
$d.Content.Find.Execute("57-28=29", $true, $false, $false, $false, $false, $false, 1, $false, "92-29=63", 2) | Out-Null
$d.Content.Find.Execute("15-4=11", $true, $false, $false, $false, $false, $false, 1, $false, "13-0=13", 2) | Out-Null
$d.Content.Find.Execute("38+5=43", $true, $false, $false, $false, $false, $false, 1, $false, "91-18=73", 2) | Out-Null
$d.Content.Find.Execute("62+23=85", $true, $false, $false, $false, $false, $false, 1, $false, "57+17=74", 2) | Out-Null
$d.Content.Find.Execute("54-39=15", $true, $false, $false, $false, $false, $false, 1, $false, "39-20=19", 2) | Out-Null
$d.Content.Find.Execute("38+23=61", $true, $false, $false, $false, $false, $false, 1, $false, "62-27=35", 2) | Out-Null
$d.Content.Find.Execute("84-69=15", $true, $false, $false, $false, $false, $false, 1, $false, "16+50=66", 2) | Out-Null
$d.Content.Find.Execute("76-28=48", $true, $false, $false, $false, $false, $false, 1, $false, "11+70=81", 2) | Out-Null
$d.Content.Find.Execute("10+83=93", $true, $false, $false, $false, $false, $false, 1, $false, "22+22=44", 2) | Out-Null
$d.Content.Find.Execute("33-11=22", $true, $false, $false, $false, $false, $false, 1, $false, "49+49=98", 2) | Out-Null
$d.Content.Find.Execute("43-5=38", $true, $false, $false, $false, $false, $false, 1, $false, "81-80=1", 2) | Out-Null
$d.Content.Find.Execute("54+34=88", $true, $false, $false, $false, $false, $false, 1, $false, "61-0=61", 2) | Out-Null
$d.Content.Find.Execute("34+59=93", $true, $false, $false, $false, $false, $false, 1, $false, "24+26=50", 2) | Out-Null
$d.Content.Find.Execute("31+68=99", $true, $false, $false, $false, $false, $false, 1, $false, "86-23=63", 2) | Out-Null
$d.Content.Find.Execute("44+53=97", $true, $false, $false, $false, $false, $false, 1, $false, "91-30=61", 2) | Out-Null
$d.Content.Find.Execute("65+30=95", $true, $false, $false, $false, $false, $false, 1, $false, "72-66=6", 2) | Out-Null
$d.Content.Find.Execute("13+62=75", $true, $false, $false, $false, $false, $false, 1, $false, "22-5=17", 2) | Out-Null
$d.Content.Find.Execute("64-2=62", $true, $false, $false, $false, $false, $false, 1, $false, "1+9=10", 2) | Out-Null
